{"js": "// Update the 25 division problems in the practice table with new values,\n// per the commit's regenerated problem set. Each old value is unique in\n// the document, so a simple search + replace per pair is safe.\nconst replacements = [\n  [\"152\u00f72=\", \"319\u00f73=\"],\n  [\"342\u00f79=\", \"460\u00f74=\"],\n  [\"775\u00f72=\", \"200\u00f77=\"],\n  [\"928\u00f77=\", \"807\u00f74=\"],\n  [\"874\u00f78=\", \"182\u00f72=\"],\n  [\"713\u00f79=\", \"889\u00f75=\"],\n  [\"712\u00f73=\", \"547\u00f73=\"],\n  [\"542\u00f76=\", \"908\u00f76=\"],\n  [\"108\u00f72=\", \"585\u00f72=\"],\n  [\"896\u00f79=\", \"506\u00f75=\"],\n  [\"334\u00f73=\", \"497\u00f72=\"],\n  [\"587\u00f78=\", \"195\u00f78=\"],\n  [\"540\u00f77=\", \"991\u00f77=\"],\n  [\"721\u00f76=\", \"612\u00f73=\"],\n  [\"742\u00f74=\", \"989\u00f72=\"],\n  [\"180\u00f75=\", \"157\u00f75=\"],\n  [\"587\u00f74=\", \"899\u00f77=\"],\n  [\"342\u00f74=\", \"858\u00f77=\"],\n  [\"478\u00f79=\", \"370\u00f78=\"],\n  [\"248\u00f74=\", \"572\u00f76=\"],\n  [\"403\u00f77=\", \"789\u00f72=\"],\n  [\"766\u00f72=\", \"964\u00f77=\"],\n  [\"550\u00f79=\", \"613\u00f79=\"],\n  [\"298\u00f76=\", \"680\u00f77=\"],\n  [\"534\u00f77=\", \"999\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 division problems in the practice table with new values,\n# per the commit's regenerated problem set. Each old value is unique in\n# the document, so a simple Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"152\u00f72=\", \"319\u00f73=\"),\n    @(\"342\u00f79=\", \"460\u00f74=\"),\n    @(\"775\u00f72=\", \"200\u00f77=\"),\n    @(\"928\u00f77=\", \"807\u00f74=\"),\n    @(\"874\u00f78=\", \"182\u00f72=\"),\n    @(\"713\u00f79=\", \"889\u00f75=\"),\n    @(\"712\u00f73=\", \"547\u00f73=\"),\n    @(\"542\u00f76=\", \"908\u00f76=\"),\n    @(\"108\u00f72=\", \"585\u00f72=\"),\n    @(\"896\u00f79=\", \"506\u00f75=\"),\n    @(\"334\u00f73=\", \"497\u00f72=\"),\n    @(\"587\u00f78=\", \"195\u00f78=\"),\n    @(\"540\u00f77=\", \"991\u00f77=\"),\n    @(\"721\u00f76=\", \"612\u00f73=\"),\n    @(\"742\u00f74=\", \"989\u00f72=\"),\n    @(\"180\u00f75=\", \"157\u00f75=\"),\n    @(\"587\u00f74=\", \"899\u00f77=\"),\n    @(\"342\u00f74=\", \"858\u00f77=\"),\n    @(\"478\u00f79=\", \"370\u00f78=\"),\n    @(\"248\u00f74=\", \"572\u00f76=\"),\n    @(\"403\u00f77=\", \"789\u00f72=\"),\n    @(\"766\u00f72=\", \"964\u00f77=\"),\n    @(\"550\u00f79=\", \"613\u00f79=\"),\n    @(\"298\u00f76=\", \"680\u00f77=\"),\n    @(\"534\u00f77=\", \"999\u00f78=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}"}
